# Insert a new data row at row 102 (pushing existing rows 102..170 down to 103..171)
# and populate it with the new record's values, matching the rest of the table's
# constant columns (A,B,C,E,F,G,H,N,Q,R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank range at row 102 (A102:R102), shifting rows 102-170 down to
# 103-171 while keeping the worksheet's used range tight (A1:R171).
$ws.Range("A102:R102").Insert(-4121)  # xlShiftDown

# Copy the row-level cell formatting from the row above (101) into the new row 102
# so that the date cell (D) keeps the date number format used throughout the table.
$ws.Range("A101:R101").Copy()
$ws.Range("A102:R102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 102 with its values
$ws.Cells.Item(102, 1).Value = 11
$ws.Cells.Item(102, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(102, 3).Value = "Bíobío"
$ws.Cells.Item(102, 4).Value = "2021-11-18"
$ws.Cells.Item(102, 5).Value = 8
$ws.Cells.Item(102, 6).Value = 100112008
$ws.Cells.Item(102, 7).Value = "Coliflor"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 1800
$ws.Cells.Item(102, 11).Value = 700
$ws.Cells.Item(102, 12).Value = 800
$ws.Cells.Item(102, 13).Value = 744
$ws.Cells.Item(102, 14).Value = "`$/unidad"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 16).Value = 744
$ws.Cells.Item(102, 17).Value = 1
$ws.Cells.Item(102, 18).Value = "Hortaliza"
